$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D4").Value = "id=`"Par1`">WHO declared SARS-CoV-2 a global pandemic.`n`n The present aim was to propose an hypothesis that there is a potential association between mean levels of vitamin D in various countries with cases and mortality caused by COVID-19. The mean levels of vitamin D for 20 European countries and morbidity and mortality caused by COVID-19 were acquired.`n`n Negative correlations between mean levels of vitamin D (average 56$([char]0x00A0)mmol/L, STDEV 10.61) in each country and the number of COVID-19 cases/1$([char]0x00A0)M (mean 295.95, STDEV 298.7, and mortality/1$([char]0x00A0)M (mean 5.96, STDEV 15.13) were observed.`n`n Vitamin D levels are severely low in the aging population especially in Spain, Italy and Switzerland.`n`n This is also the most vulnerable group of the population in relation to COVID-19. It should be advisable to perform dedicated studies about vitamin D levels in COVID-19 patients with different degrees of disease severity.`n`n"
$ws.Range("E4").Value = "[Petre Cristian%Ilie%dr.cristianilie@gmail.com%1,  Simina%Stefanescu%NULL%1,  Lee%Smith%NULL%1]"

# Restore row 4's default (non-custom) height, since assigning the long
# abstract text can trigger an implicit auto row-height flag.
$ws.Rows.Item(4).AutoFit()
